# Added crdc login backup code:
# - The backup code previously sitting at A11 becomes the active/next code in A2.
# - The consumed codes (old A2, A3, A4) and the now-vacated A11 are removed.
# - A12 (the last remaining backup code) is left untouched.
# - Selection moves to the new active code cell, A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Promote the code that was in A11 into A2 (the "current" backup code slot).
$ws.Range("A2").Value = $ws.Range("A11").Value2

# Remove the now-used-up codes and the vacated row.
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("A11").ClearContents()

# Match the saved selection state.
$ws.Range("A2").Select()
